$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

Set-TextValue $ws.Range("D2") '276.39'
Set-TextValue $ws.Range("E2") '1.46%'
Set-TextValue $ws.Range("G2") '2'

Set-TextValue $ws.Range("D3") '27.29'
Set-TextValue $ws.Range("E3") '1.91%'
Set-TextValue $ws.Range("G3") '2'

Set-TextValue $ws.Range("E4") '-0.70%'
Set-TextValue $ws.Range("G4") '2'

Set-TextValue $ws.Range("D5") '0.06411'
Set-TextValue $ws.Range("E5") '1.25%'
Set-TextValue $ws.Range("G5") '2'

Set-TextValue $ws.Range("D6") '6.931'
Set-TextValue $ws.Range("E6") '0.65%'
Set-TextValue $ws.Range("G6") '2'

Set-TextValue $ws.Range("D7") '1.210'
Set-TextValue $ws.Range("E7") '-7.67%'
Set-TextValue $ws.Range("G7") '2'

Set-TextValue $ws.Range("D8") '0.8811'
Set-TextValue $ws.Range("E8") '0.04%'
Set-TextValue $ws.Range("G8") '2'

Set-TextValue $ws.Range("D9") '0.1516'
Set-TextValue $ws.Range("E9") '3.88%'
Set-TextValue $ws.Range("G9") '2'

Set-TextValue $ws.Range("D10") '0.05070'
Set-TextValue $ws.Range("E10") '-0.30%'
Set-TextValue $ws.Range("G10") '2'

Set-TextValue $ws.Range("D11") '0.07593'
Set-TextValue $ws.Range("E11") '3.38%'
Set-TextValue $ws.Range("G11") '2'

Set-TextValue $ws.Range("D12") '0.02963'
Set-TextValue $ws.Range("E12") '-4.55%'
Set-TextValue $ws.Range("G12") '2'

Set-TextValue $ws.Range("D13") '0.09007'
Set-TextValue $ws.Range("E13") '-0.37%'
Set-TextValue $ws.Range("G13") '2'

Set-TextValue $ws.Range("D14") '0.001565'
Set-TextValue $ws.Range("E14") '-0.38%'
Set-TextValue $ws.Range("G14") '2'

Set-TextValue $ws.Range("D15") '0.0006424'
Set-TextValue $ws.Range("E15") '1.48%'
Set-TextValue $ws.Range("G15") '2'

Set-TextValue $ws.Range("D16") '0.006187'
Set-TextValue $ws.Range("E16") '2.50%'
Set-TextValue $ws.Range("G16") '2'

Set-TextValue $ws.Range("D17") '3.469'
Set-TextValue $ws.Range("E17") '0.05%'
Set-TextValue $ws.Range("G17") '2'

Set-TextValue $ws.Range("D18") '3.304'
Set-TextValue $ws.Range("E18") '-1.44%'
Set-TextValue $ws.Range("G18") '2'

Set-TextValue $ws.Range("D19") '2.285'
Set-TextValue $ws.Range("E19") '0.05%'
Set-TextValue $ws.Range("G19") '2'

Set-TextValue $ws.Range("G20") '2'

Set-TextValue $ws.Range("D21") '0.1355'
Set-TextValue $ws.Range("E21") '2.10%'
Set-TextValue $ws.Range("G21") '2'

Set-TextValue $ws.Range("D22") '3.918'
Set-TextValue $ws.Range("E22") '-0.09%'
Set-TextValue $ws.Range("G22") '2'

Set-TextValue $ws.Range("D23") '0.04430'
Set-TextValue $ws.Range("E23") '-0.08%'
Set-TextValue $ws.Range("G23") '2'

Set-TextValue $ws.Range("D24") '0.001175'
Set-TextValue $ws.Range("E24") '-0.28%'
Set-TextValue $ws.Range("G24") '2'

Set-TextValue $ws.Range("D25") '0.004266'
Set-TextValue $ws.Range("E25") '15.58%'
Set-TextValue $ws.Range("G25") '2'

Set-TextValue $ws.Range("D26") '0.0001201'
Set-TextValue $ws.Range("E26") '-0.33%'
Set-TextValue $ws.Range("G26") '2'

Set-TextValue $ws.Range("D27") '0.0001939'
Set-TextValue $ws.Range("E27") '13.78%'
Set-TextValue $ws.Range("G27") '2'

Set-TextValue $ws.Range("G28") '2'

Set-TextValue $ws.Range("G29") '2'

Set-TextValue $ws.Range("G30") '2'

Set-TextValue $ws.Range("G31") '2'

Set-TextValue $ws.Range("G32") '2'

Set-TextValue $ws.Range("G33") '2'

Set-TextValue $ws.Range("G34") '2'

Set-TextValue $ws.Range("G35") '2'

Set-TextValue $ws.Range("G36") '2'

Set-TextValue $ws.Range("G37") '2'

Set-TextValue $ws.Range("G38") '2'

Set-TextValue $ws.Range("G39") '2'

Set-TextValue $ws.Range("D40") '0.04158'
Set-TextValue $ws.Range("E40") '2.26%'
Set-TextValue $ws.Range("G40") '2'

Set-TextValue $ws.Range("D41") '0.006829'
Set-TextValue $ws.Range("E41") '3.10%'
Set-TextValue $ws.Range("G41") '2'

Set-TextValue $ws.Range("D42") '0.1175'
Set-TextValue $ws.Range("E42") '0.86%'
Set-TextValue $ws.Range("G42") '2'

Set-TextValue $ws.Range("D43") '0.002132'
Set-TextValue $ws.Range("E43") '0.62%'
Set-TextValue $ws.Range("G43") '2'

Set-TextValue $ws.Range("D44") '0.01186'
Set-TextValue $ws.Range("E44") '-1.17%'
Set-TextValue $ws.Range("G44") '2'

Set-TextValue $ws.Range("D45") '0.00005175'
Set-TextValue $ws.Range("E45") '-2.70%'
Set-TextValue $ws.Range("G45") '2'

Set-TextValue $ws.Range("B46") 'BOLO'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws.Range("D46") '1.510'
Set-TextValue $ws.Range("E46") '-35.92%'
Set-TextValue $ws.Range("G46") '2'

Set-TextValue $ws.Range("B47") 'CoinbaseStockToken'
Set-TextValue $ws.Range("C47") 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws.Range("D47") '0.02002'
Set-TextValue $ws.Range("E47") '-2.78%'
Set-TextValue $ws.Range("G47") '2'

Set-TextValue $ws.Range("G48") '2'

Set-TextValue $ws.Range("G49") '2'

Set-TextValue $ws.Range("G50") '2'

Set-TextValue $ws.Range("G51") '2'
